$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$startLen = $tr.Length

# Append the two new paragraphs as plain text first (paragraph marks via `r).
$tr.InsertAfter("`rВывод:`rИгра реализована качественно и имеет возможности для улучшения. Удовлетворяет своим целям.")

$full = $shp.TextFrame.TextRange

$p1Start = $startLen + 2
$word1 = $full.Characters($p1Start, 5)
$word1.Text = "Вывод"
$word1.Font.Italic = $false
$word1.LanguageID = 1049

$colon = $full.Characters($p1Start + 5, 1)
$colon.Text = ":"
$colon.Font.Italic = $false
$colon.LanguageID = 1033
